# Refresh the crypto-symbol table: updated Price/Volume(1h) figures, plus
# a 3-way row rotation among BKEXToken / CEJI / KickToken (rows 41-43).
# NumberFormat is forced to Text ("@") before each Price/Volume write so the
# numeric-looking strings ("245.36", "-0.74%", ...) stay text values instead
# of being auto-coerced into numbers/percentages by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.36"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.74%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.108"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.28%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05688"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.56%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.526"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.60%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8194"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.80%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8598"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.81%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.42%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06947"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.81%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02860"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.77%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09390"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.07%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001540"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.58%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04068"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-12.51%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005980"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-94.00%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006215"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.37%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.507"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.77%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.32%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.317"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "12.73%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3165"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03222"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.35%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.08%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.552"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.30%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.79%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001216"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.25%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004470"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001180"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "22.90%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001445"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-25.41%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03723"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.60%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005978"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.78%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1058"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.09%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002299"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.01%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009715"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "17.87%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005106"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.32%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.01%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-30.33%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.95%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
